$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3, shifting the old rows 3-8 (and the
# trailing spacer rows) down by one.
$ws.Rows("3:3").Insert()

# --- New row 3 (inserted): "rdp = Shell_RDP(); util = EasyshellLib.CommonUtils" ---
$ws.Range("A3").Value = "N"
$ws.Range("B3").Value = "rdp = Shell_RDP(); util = EasyshellLib.CommonUtils"

# Give B3 the same cell style as B4 (header-ish style s=7) instead of the
# plain style Insert() copied down from row 2.
$ws.Range("B4").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Match the row height used by the rest of this block of rows.
$ws.Rows("3:3").RowHeight = 15.75

# --- Row 4 (was row 3): Shell_RDP().create('standardRDP') -> rdp.create('standardRDP') ---
$cell4 = $ws.Range("B4")
$cell4.Value = "rdp.create('standardRDP')"
$cell4.Characters(12, 13).Font.Bold = $true
$cell4.Characters(12, 13).Font.Color = 8421376
$cell4.Characters(25, 1).Font.Bold = $false
$cell4.Characters(25, 1).Font.Color = 0

# --- Row 5 (was row 4): SwitchToUser ---
$ws.Range("A5").Value = "#N"
$ws.Range("B5").Value = "util.SwitchToUser()"

# --- Row 6 (was row 5): Reboot ---
$ws.Range("A6").Value = "#N"
$ws.Range("B6").Value = "util.Reboot()"

# --- Row 7 (was row 6): check('standardRDP') ---
$ws.Range("A7").Value = "#Y"
$ws.Range("B7").Value = "rdp.check('standardRDP')"

# --- Row 8 (was row 7): SwitchToAdmin ---
$ws.Range("A8").Value = "#N"
$ws.Range("B8").Value = "util.SwitchToAdmin()"

# --- Row 9 (was row 8): Reboot ---
$ws.Range("A9").Value = "#N"
$ws.Range("B9").Value = "util.Reboot()"

# Selection moves to B15
$ws.Range("B15").Select()
